#
# Adds "Errors" and "Warnings" sheets after "Classes", populates the
# "Warnings" sheet with three messages, and makes "Warnings" the active tab.
#
# The messages all start with a literal apostrophe ('). Typing/assigning a
# string that begins with an apostrophe directly into a cell is interpreted
# by Excel (and this host) as the "treat as text" quote-prefix marker: the
# apostrophe is stripped from the stored value and a quotePrefix style is
# applied instead. To get a *literal* leading apostrophe in the stored
# shared-string text (with a plain, unstyled cell) we build the text with a
# formula (="'" & "...") and then convert that formula to a static value via
# Copy / PasteSpecial-values, exactly as a user would do in the Excel UI.

$wb = $excel.ActiveWorkbook
$wsClasses = $wb.Worksheets.Item(1)

# New sheets, inserted right after "Classes".
$wsErrors = $wb.Worksheets.Add($null, $wsClasses)
$wsErrors.Name = "Errors"

$wsWarnings = $wb.Worksheets.Add($null, $wsErrors)
$wsWarnings.Name = "Warnings"

# The three warning messages (row numbers from the original import).
$rows = 4, 7, 9

for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowNum = $rows[$i]
    $cell = $wsWarnings.Cells.Item($i + 1, 1)

    $text = 'Sheet ""Classes"" Row: ' + $rowNum + ' No data found between cells ""A"" and ""D"" Skipping this row'','
    $formula = '="' + "'" + '"&"' + $text + '"'

    $cell.Formula = $formula
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = $false

# Match the original sheet's selection on "Warnings" and make it the active
# (selected) sheet/tab.
$wsWarnings.Range("B7:B8").Select()

# "Classes" no longer keeps the Excel-UI selected tab (Warnings becomes it);
# restore its own cell selection as recorded before the edit.
$wsClasses.Range("A7:XFD7").Select()
$wsWarnings.Activate()
$wsWarnings.Range("B7:B8").Select()
